$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 18250
$ws.Range("J32").Value = 18250
$ws.Range("L32").Value = 18250
$ws.Range("N32").Value = -18902

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 8132.8
$ws.Range("I43").Value = 2750
$ws.Range("K43").Value = 2750
$ws.Range("M43").Value = -2681

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 6024
$ws.Range("I53").Value = 1036.2727
$ws.Range("J53").Value = 11510.5
$ws.Range("K53").Value = 1036.2727
$ws.Range("L53").Value = 11510.5
$ws.Range("M53").Value = -399.2727
$ws.Range("N53").Value = -12784.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 4031.4358
$ws.Range("J112").Value = 2421.697
$ws.Range("L112").Value = 7265.091
$ws.Range("N112").Value = -9481.091

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 21907.666
$ws.Range("I113").Value = 22564.416
$ws.Range("J113").Value = 21250.916
$ws.Range("K113").Value = 22564.416
$ws.Range("L113").Value = 21250.916
$ws.Range("M113").Value = -19310.416
$ws.Range("N113").Value = -27758.916

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 4634.451
$ws.Range("I132").Value = 3405.6934
$ws.Range("K132").Value = 10217.0802
$ws.Range("M132").Value = -7687.0802

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 6160.56
$ws.Range("I138").Value = 2519.9092
$ws.Range("J138").Value = 7953.7163
$ws.Range("K138").Value = 7559.7276
$ws.Range("L138").Value = 23861.1489
$ws.Range("M138").Value = -2419.7276
$ws.Range("N138").Value = -34141.1489

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3079.8823
$ws.Range("I45").Value = 3240.8572
$ws.Range("K45").Value = 3240.8572
$ws.Range("M45").Value = -2863.8572

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 23100.857
$ws.Range("I74").Value = 3927.4285
$ws.Range("J74").Value = 32687.572
$ws.Range("K74").Value = 3927.4285
$ws.Range("L74").Value = 32687.572
$ws.Range("M74").Value = -3053.4285
$ws.Range("N74").Value = -34435.572

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 23100.857
$ws.Range("I77").Value = 3927.4285
$ws.Range("J77").Value = 32687.572
$ws.Range("K77").Value = 19637.1425
$ws.Range("L77").Value = 163437.86
$ws.Range("M77").Value = -15269.1425
$ws.Range("N77").Value = -172173.86

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 17768.021
$ws.Range("I20").Value = 4378.96
$ws.Range("J20").Value = 33707.383
$ws.Range("K20").Value = 4378.96
$ws.Range("L20").Value = 33707.383
$ws.Range("M20").Value = -4131.96
$ws.Range("N20").Value = -34201.383

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H88").Value = 32900
$ws.Range("J88").Value = 33125
$ws.Range("L88").Value = 33125
$ws.Range("N88").Value = -33937

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H91").Value = 32900
$ws.Range("J91").Value = 33125
$ws.Range("L91").Value = 33125
$ws.Range("N91").Value = -35933

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 6072.3335
$ws.Range("I94").Value = 3160.4614
$ws.Range("K94").Value = 3160.4614
$ws.Range("M94").Value = -2709.4614

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3056.4285
$ws.Range("I105").Value = 3056.4285
$ws.Range("K105").Value = 3056.4285
$ws.Range("M105").Value = -1309.4285

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 8000.933
$ws.Range("I134").Value = 2516.1562
$ws.Range("K134").Value = 7548.4686
$ws.Range("M134").Value = -5013.4686

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 20809.648
$ws.Range("I31").Value = 6514.6665
$ws.Range("K31").Value = 6514.6665
$ws.Range("M31").Value = -6219.6665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 20809.648
$ws.Range("I34").Value = 6514.6665
$ws.Range("K34").Value = 6514.6665
$ws.Range("M34").Value = -6312.6665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 14152.559
$ws.Range("J58").Value = 27997.46
$ws.Range("L58").Value = 27997.46
$ws.Range("N58").Value = -28403.46

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 5628.1
$ws.Range("I62").Value = 5879.1665
$ws.Range("K62").Value = 5879.1665
$ws.Range("M62").Value = -5255.1665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 5628.1
$ws.Range("I65").Value = 5879.1665
$ws.Range("K65").Value = 29395.8325
$ws.Range("M65").Value = -26275.8325

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 5021.0454
$ws.Range("I86").Value = 4278.8184
$ws.Range("K86").Value = 4278.8184
$ws.Range("M86").Value = -3155.8184

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 5021.0454
$ws.Range("I89").Value = 4278.8184
$ws.Range("K89").Value = 21394.092
$ws.Range("M89").Value = -15778.092

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H120").Value = 20833.262
$ws.Range("J120").Value = 20833.262
$ws.Range("L120").Value = 20833.262
$ws.Range("N120").Value = -28091.262

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 14152.559
$ws.Range("J136").Value = 27997.46
$ws.Range("L136").Value = 83992.38
$ws.Range("N136").Value = -89092.38

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 235.38461
$ws.Range("I2").Value = 260.27777
$ws.Range("J2").Value = 179.375
$ws.Range("K2").Value = 1561.66662
$ws.Range("L2").Value = 1076.25
$ws.Range("M2").Value = -1448.66662
$ws.Range("N2").Value = -1302.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 3176631.5
$ws.Range("I34").Value = 1133.1666
$ws.Range("J34").Value = 7410629
$ws.Range("K34").Value = 3399.4998
$ws.Range("L34").Value = 22231887
$ws.Range("M34").Value = -3315.4998
$ws.Range("N34").Value = -22232055

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H45").Value = 2464.6667
$ws.Range("I45").Value = 1197
$ws.Range("J45").Value = 5000
$ws.Range("K45").Value = 3591
$ws.Range("L45").Value = 15000
$ws.Range("M45").Value = -3059
$ws.Range("N45").Value = -16064

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 4522
$ws.Range("I46").Value = 518
$ws.Range("J46").Value = 7525
$ws.Range("K46").Value = 1554
$ws.Range("L46").Value = 22575
$ws.Range("M46").Value = -1463
$ws.Range("N46").Value = -22757

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 11595.462
$ws.Range("I121").Value = 549.8
$ws.Range("J121").Value = 18499
$ws.Range("K121").Value = 1649.4
$ws.Range("L121").Value = 55497
$ws.Range("M121").Value = -339.3999999999999
$ws.Range("N121").Value = -58117

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 4498.8
$ws.Range("J137").Value = 4999.6665
$ws.Range("L137").Value = 14998.9995
$ws.Range("N137").Value = -25198.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 5572.2104
$ws.Range("I139").Value = 4983.857
$ws.Range("K139").Value = 14951.571
$ws.Range("M139").Value = -9811.571

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 1184.5454
$ws.Range("I140").Value = 1184.5454
$ws.Range("K140").Value = 3553.6362
$ws.Range("M140").Value = 1626.3638

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 21582.25
$ws.Range("J80").Value = 28999
$ws.Range("L80").Value = 28999
$ws.Range("N80").Value = -30995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 21582.25
$ws.Range("J83").Value = 28999
$ws.Range("L83").Value = 144995
$ws.Range("N83").Value = -154979

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3966.182
$ws.Range("I102").Value = 5683.5
$ws.Range("J102").Value = 1905.4
$ws.Range("K102").Value = 5683.5
$ws.Range("L102").Value = 1905.4
$ws.Range("M102").Value = -4061.5
$ws.Range("N102").Value = -5149.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 7515.9614
$ws.Range("I132").Value = 3695.5454
$ws.Range("K132").Value = 11086.6362
$ws.Range("M132").Value = -8556.636200000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 7820.0933
$ws.Range("I132").Value = 4667.3076
$ws.Range("K132").Value = 14001.9228
$ws.Range("M132").Value = -11471.9228

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 23299.8
$ws.Range("J31").Value = 23299.8
$ws.Range("L31").Value = 23299.8
$ws.Range("N31").Value = -23995.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 1923.0769
$ws.Range("I54").Value = 1000
$ws.Range("J54").Value = 3000
$ws.Range("K54").Value = 1000
$ws.Range("L54").Value = 3000
$ws.Range("M54").Value = -480
$ws.Range("N54").Value = -4040

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5120.875
$ws.Range("J62").Value = 5666.6665
$ws.Range("L62").Value = 5666.6665
$ws.Range("N62").Value = -6914.6665

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 5120.875
$ws.Range("J65").Value = 5666.6665
$ws.Range("L65").Value = 28333.3325
$ws.Range("N65").Value = -34573.3325

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2937.6365
$ws.Range("I122").Value = 2033.6428
$ws.Range("J122").Value = 8000
$ws.Range("K122").Value = 6100.928400000001
$ws.Range("L122").Value = 24000
$ws.Range("M122").Value = -3650.928400000001
$ws.Range("N122").Value = -28900

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 9617.587
$ws.Range("I132").Value = 5232.484
$ws.Range("K132").Value = 15697.452
$ws.Range("M132").Value = -13167.452

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 14207.389
$ws.Range("I136").Value = 2152.5
$ws.Range("K136").Value = 6457.5
$ws.Range("M136").Value = -3907.5
